$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Add two new authoring test case rows (E40, E41) -----------------------
# The write order below matters: it reproduces the exact shared-string
# first-use order recorded in the source workbook (the "post" description
# was typed into row 42 before the "patent" description was typed into
# row 41, even though row 41 precedes row 42 on the sheet).

$ws.Range("A41").Value = "TestCase_E40"
$ws.Range("B41").Value = "OPQA-1108"

$ws.Range("C42").Value = "Verify that same post can be added to multiple watchlists"
$ws.Range("C41").Value = "Verify that same patent can be added to multiple watchlists"

$ws.Range("A42").Value = "TestCase_E41"
$ws.Range("B42").Value = "OPQA-1109"

$ws.Range("D41").Value = "Y"
$ws.Range("E41").Value = "PASS"
$ws.Range("D42").Value = "Y"
$ws.Range("E42").Value = "PASS"

# Match the formatting (borders etc.) of the preceding data row for the two
# newly added rows, the same way a user would fill-down / copy formatting.
[void]$ws.Range("A40:E40").Copy()
[void]$ws.Range("A41:E42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- View state: last active selection lands on the new description cell --
[void]$ws.Range("C41").Select()

# --- Window size, as recorded after the edit session ------------------------
$win = $excel.ActiveWindow
$win.Width = 12240
$win.Height = 10125
